# Generate Report for Handback
# Updates the localization-status workbook to reflect that the two handed-off
# files now have a "handed back" / target / handback-file / handback-datetime
# recorded, and widens a few columns that now need to show the longer status
# text / hyperlinks.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    on every sheet that shows it.
# ---------------------------------------------------------------------------
$ws1.Range("E2").Value2 = $newStatus
$ws1.Range("F2").Value2 = $newStatus
$ws1.Range("E3").Value2 = $newStatus
$ws1.Range("F3").Value2 = $newStatus

$ws2.Range("C2").Value2 = $newStatus
$ws2.Range("C3").Value2 = $newStatus

$ws3.Range("C2").Value2 = $newStatus
$ws3.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn / de-de sheets: fill in "Latest Target File" (I) and
#    "Latest Handback File" (J) columns for both data rows, and turn the
#    target-file cell into a hyperlink to the same markdown file as column A.
# ---------------------------------------------------------------------------
$urlMd1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/6364f5345d6395a7c0fac3f44dfc36db00b759dc/e2e/0c7147df-1fe2-4255-a57c-9730c6e48ba4.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/6364f5345d6395a7c0fac3f44dfc36db00b759dc/e2e/1b5f9af8-8438-40ec-95f4-0b7e461e3c53.md"

foreach ($ws in @($ws2, $ws3)) {
    # Add the new hyperlinks for row 2 first...
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlMd1, "", "", "0c7147df-1fe2-4255-a57c-9730c6e48ba4.md")

    # ...then recreate the existing row-3 hyperlink so the relationship /
    # hyperlink ordering interleaves as A2, I2, A3, I3 (matching the way the
    # report generator lays the sheet out).
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq "`$A`$3") {
            $hl.Delete()
        }
    }
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlMd2, "", "", "1b5f9af8-8438-40ec-95f4-0b7e461e3c53.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlMd2, "", "", "1b5f9af8-8438-40ec-95f4-0b7e461e3c53.md")

    # Match the style used by the existing A2/A3 hyperlink cells.
    $ws.Range("I2").Style = "HyperLink"
    $ws.Range("I3").Style = "HyperLink"
    $ws.Range("A3").Style = "HyperLink"
}

# Latest Target File display text (same file name as the handoff file).
$ws2.Range("I2").Value2 = "0c7147df-1fe2-4255-a57c-9730c6e48ba4.md"
$ws2.Range("I3").Value2 = "1b5f9af8-8438-40ec-95f4-0b7e461e3c53.md"
$ws3.Range("I2").Value2 = "0c7147df-1fe2-4255-a57c-9730c6e48ba4.md"
$ws3.Range("I3").Value2 = "1b5f9af8-8438-40ec-95f4-0b7e461e3c53.md"

# Latest Handback File: generated xliff file per language.
$ws2.Range("J2").Value2 = "0c7147df-1fe2-4255-a57c-9730c6e48ba4.c3a034f60ab2c4b0aba14bbd806fb56c3a22db62.zh-cn.xlf"
$ws2.Range("J3").Value2 = "1b5f9af8-8438-40ec-95f4-0b7e461e3c53.d3179f39da597e559fa0030b66c031617e576e07.zh-cn.xlf"
$ws3.Range("J2").Value2 = "0c7147df-1fe2-4255-a57c-9730c6e48ba4.c3a034f60ab2c4b0aba14bbd806fb56c3a22db62.de-de.xlf"
$ws3.Range("J3").Value2 = "1b5f9af8-8438-40ec-95f4-0b7e461e3c53.d3179f39da597e559fa0030b66c031617e576e07.de-de.xlf"

# Latest Handback DateTime: zh-cn rows share one timestamp (already showed
# the placeholder 0001-01-01 00:00:00, so simply overwrite it); de-de rows
# get their own, later timestamp.
$ws2.Range("K2").Value2 = "2016-08-12 14:32:04"
$ws2.Range("K3").Value2 = "2016-08-12 14:32:04"
$ws3.Range("K2").Value2 = "2016-08-12 14:32:18"
$ws3.Range("K3").Value2 = "2016-08-12 14:32:18"

# ---------------------------------------------------------------------------
# 3. Widen columns that now hold longer hyperlink / status text.
# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) -> same width as the
# "Latest HO Xliff Generate Date" column (G).
$ws1.Columns.Item(5).ColumnWidth = 29.17
$ws1.Columns.Item(6).ColumnWidth = 29.17

foreach ($ws in @($ws2, $ws3)) {
    # Status column (C) widens the same way as the Overview columns.
    $ws.Columns.Item(3).ColumnWidth = 29.17
    # Latest Target File (I) / Latest Handback File (J) now hold hyperlinks
    # and long generated file names, so match the other wide (40-char)
    # columns already on the sheet.
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}
